{"js": "// Replace the date line and each \"a op b=\" cell text with its new value.\n// Every \"old\" value is unique within the document, so an exact,\n// case-sensitive, whole-string search reliably targets the single\n// matching run in each case.\nconst replacements = [\n  [\"2024-08-09 Friday\", \"2024-08-10 Saturday\"],\n  [\"94-44=\", \"78-9=\"],\n  [\"61-58=\", \"96-88=\"],\n  [\"14-11=\", \"48-6=\"],\n  [\"39+60=\", \"25+53=\"],\n  [\"3+20=\", \"88-3=\"],\n  [\"9+39=\", \"41-25=\"],\n  [\"20+18=\", \"80-33=\"],\n  [\"6+67=\", \"99-40=\"],\n  [\"9+52=\", \"9+49=\"],\n  [\"71-34=\", \"52-47=\"],\n  [\"49+44=\", \"78-34=\"],\n  [\"23+45=\", \"57-14=\"],\n  [\"61-28=\", \"42-0=\"],\n  [\"32-0=\", \"7+9=\"],\n  [\"23+44=\", \"81+8=\"],\n  [\"57-38=\", \"27+19=\"],\n  [\"76+3=\", \"96-25=\"],\n  [\"55+37=\", \"44-8=\"],\n  [\"67+21=\", \"73-56=\"],\n  [\"47+19=\", \"94-93=\"],\n  [\"56-22=\", \"18+30=\"],\n  [\"88-24=\", \"65-65=\"],\n  [\"71+20=\", \"78-3=\"],\n  [\"17+77=\", \"42+57=\"],\n  [\"92-63=\", \"68-47=\"],\n  [\"20+74=\", \"85-63=\"],\n  [\"36-29=\", \"15+67=\"],\n  [\"49+23=\", \"49-30=\"],\n  [\"56+27=\", \"11-8=\"],\n  [\"57+33=\", \"98-93=\"],\n  [\"36+57=\", \"60+27=\"],\n  [\"59+36=\", \"39+1=\"],\n  [\"82-15=\", \"33-4=\"],\n  [\"67+12=\", \"69+8=\"],\n  [\"84-65=\", \"33+8=\"],\n  [\"95-2=\", \"74-1=\"],\n  [\"13+57=\", \"48+0=\"],\n  [\"94-86=\", \"13+68=\"],\n  [\"62+36=\", \"17+81=\"],\n  [\"34-1=\", \"65+15=\"],\n  [\"65+6=\", \"59+13=\"],\n  [\"86-56=\", \"31+26=\"],\n  [\"14-7=\", \"56-0=\"],\n  [\"59-12=\", \"64+2=\"],\n  [\"4+16=\", \"28+37=\"],\n  [\"96-26=\", \"65+10=\"],\n  [\"75-37=\", \"91-87=\"],\n  [\"62-41=\", \"42+24=\"],\n  [\"54+38=\", \"28+28=\"],\n  [\"41-23=\", \"59-46=\"],\n  [\"45+43=\", \"42+52=\"],\n  [\"26+32=\", \"93-76=\"],\n  [\"82-8=\", \"75-24=\"],\n  [\"25+9=\", \"17+80=\"],\n  [\"64-52=\", \"23-0=\"],\n  [\"70-43=\", \"40+18=\"],\n  [\"62+22=\", \"28+46=\"],\n  [\"70-57=\", \"96-65=\"],\n  [\"8+41=\", \"55-30=\"],\n  [\"11+37=\", \"77+2=\"],\n  [\"36+11=\", \"23+69=\"],\n  [\"13+51=\", \"93-52=\"],\n  [\"17+50=\", \"72-29=\"],\n  [\"19+78=\", \"91-48=\"],\n  [\"90-52=\", \"24+22=\"],\n  [\"78-12=\", \"63-26=\"],\n  [\"25+43=\", \"65+31=\"],\n  [\"29-1=\", \"43+29=\"],\n  [\"53-28=\", \"89-37=\"],\n  [\"29-18=\", \"94-3=\"],\n  [\"39+7=\", \"3+34=\"],\n  [\"49-1=\", \"83-37=\"],\n  [\"9+57=\", \"64+16=\"],\n  [\"54+15=\", \"35+16=\"],\n  [\"35+47=\", \"29+18=\"],\n  [\"34+22=\", \"85-55=\"],\n  [\"1+96=\", \"84-63=\"],\n  [\"21+74=\", \"21+22=\"],\n  [\"28+51=\", \"47+8=\"],\n  [\"23+28=\", \"74+16=\"],\n  [\"38-3=\", \"47-39=\"],\n  [\"80-2=\", \"56-52=\"],\n  [\"15+56=\", \"24+25=\"],\n  [\"54+45=\", \"44-22=\"],\n  [\"1+45=\", \"28+50=\"],\n  [\"51-20=\", \"97-37=\"],\n  [\"50-16=\", \"53-43=\"],\n  [\"61-31=\", \"75-39=\"],\n  [\"98-20=\", \"29+14=\"],\n  [\"57-13=\", \"49-45=\"],\n  [\"70+15=\", \"19+74=\"],\n  [\"28+49=\", \"36-0=\"],\n  [\"29+9=\", \"53-29=\"],\n  [\"98-82=\", \"77-70=\"],\n  [\"95-94=\", \"29+30=\"],\n  [\"60-2=\", \"29+4=\"],\n  [\"13-1=\", \"17+16=\"],\n  [\"16+40=\", \"8-5=\"],\n  [\"45-22=\", \"61+15=\"],\n  [\"49+2=\", \"97-39=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"a op b=\" cell text with its new value.\n# Every \"old\" value is unique within the document, so an exact,\n# case-sensitive Find/Replace (wdReplaceAll, MatchCase=$true,\n# MatchWholeWord=$false) reliably targets only the intended text.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old='2024-08-09 Friday'; New='2024-08-10 Saturday'},\n    @{Old='94-44='; New='78-9='},\n    @{Old='61-58='; New='96-88='},\n    @{Old='14-11='; New='48-6='},\n    @{Old='39+60='; New='25+53='},\n    @{Old='3+20='; New='88-3='},\n    @{Old='9+39='; New='41-25='},\n    @{Old='20+18='; New='80-33='},\n    @{Old='6+67='; New='99-40='},\n    @{Old='9+52='; New='9+49='},\n    @{Old='71-34='; New='52-47='},\n    @{Old='49+44='; New='78-34='},\n    @{Old='23+45='; New='57-14='},\n    @{Old='61-28='; New='42-0='},\n    @{Old='32-0='; New='7+9='},\n    @{Old='23+44='; New='81+8='},\n    @{Old='57-38='; New='27+19='},\n    @{Old='76+3='; New='96-25='},\n    @{Old='55+37='; New='44-8='},\n    @{Old='67+21='; New='73-56='},\n    @{Old='47+19='; New='94-93='},\n    @{Old='56-22='; New='18+30='},\n    @{Old='88-24='; New='65-65='},\n    @{Old='71+20='; New='78-3='},\n    @{Old='17+77='; New='42+57='},\n    @{Old='92-63='; New='68-47='},\n    @{Old='20+74='; New='85-63='},\n    @{Old='36-29='; New='15+67='},\n    @{Old='49+23='; New='49-30='},\n    @{Old='56+27='; New='11-8='},\n    @{Old='57+33='; New='98-93='},\n    @{Old='36+57='; New='60+27='},\n    @{Old='59+36='; New='39+1='},\n    @{Old='82-15='; New='33-4='},\n    @{Old='67+12='; New='69+8='},\n    @{Old='84-65='; New='33+8='},\n    @{Old='95-2='; New='74-1='},\n    @{Old='13+57='; New='48+0='},\n    @{Old='94-86='; New='13+68='},\n    @{Old='62+36='; New='17+81='},\n    @{Old='34-1='; New='65+15='},\n    @{Old='65+6='; New='59+13='},\n    @{Old='86-56='; New='31+26='},\n    @{Old='14-7='; New='56-0='},\n    @{Old='59-12='; New='64+2='},\n    @{Old='4+16='; New='28+37='},\n    @{Old='96-26='; New='65+10='},\n    @{Old='75-37='; New='91-87='},\n    @{Old='62-41='; New='42+24='},\n    @{Old='54+38='; New='28+28='},\n    @{Old='41-23='; New='59-46='},\n    @{Old='45+43='; New='42+52='},\n    @{Old='26+32='; New='93-76='},\n    @{Old='82-8='; New='75-24='},\n    @{Old='25+9='; New='17+80='},\n    @{Old='64-52='; New='23-0='},\n    @{Old='70-43='; New='40+18='},\n    @{Old='62+22='; New='28+46='},\n    @{Old='70-57='; New='96-65='},\n    @{Old='8+41='; New='55-30='},\n    @{Old='11+37='; New='77+2='},\n    @{Old='36+11='; New='23+69='},\n    @{Old='13+51='; New='93-52='},\n    @{Old='17+50='; New='72-29='},\n    @{Old='19+78='; New='91-48='},\n    @{Old='90-52='; New='24+22='},\n    @{Old='78-12='; New='63-26='},\n    @{Old='25+43='; New='65+31='},\n    @{Old='29-1='; New='43+29='},\n    @{Old='53-28='; New='89-37='},\n    @{Old='29-18='; New='94-3='},\n    @{Old='39+7='; New='3+34='},\n    @{Old='49-1='; New='83-37='},\n    @{Old='9+57='; New='64+16='},\n    @{Old='54+15='; New='35+16='},\n    @{Old='35+47='; New='29+18='},\n    @{Old='34+22='; New='85-55='},\n    @{Old='1+96='; New='84-63='},\n    @{Old='21+74='; New='21+22='},\n    @{Old='28+51='; New='47+8='},\n    @{Old='23+28='; New='74+16='},\n    @{Old='38-3='; New='47-39='},\n    @{Old='80-2='; New='56-52='},\n    @{Old='15+56='; New='24+25='},\n    @{Old='54+45='; New='44-22='},\n    @{Old='1+45='; New='28+50='},\n    @{Old='51-20='; New='97-37='},\n    @{Old='50-16='; New='53-43='},\n    @{Old='61-31='; New='75-39='},\n    @{Old='98-20='; New='29+14='},\n    @{Old='57-13='; New='49-45='},\n    @{Old='70+15='; New='19+74='},\n    @{Old='28+49='; New='36-0='},\n    @{Old='29+9='; New='53-29='},\n    @{Old='98-82='; New='77-70='},\n    @{Old='95-94='; New='29+30='},\n    @{Old='60-2='; New='29+4='},\n    @{Old='13-1='; New='17+16='},\n    @{Old='16+40='; New='8-5='},\n    @{Old='45-22='; New='61+15='},\n    @{Old='49+2='; New='97-39='}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $r.Old,      # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $r.New,      # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"No match found for: $($r.Old)\"\n    }\n}\n"}
